# Generate Report for Handoff
# Adds a new file entry (d6da7312-fbe8-491d-afc1-66ef03145004) as row 7
# on the Overview, zh-cn and de-de sheets, matching a "Ready for
# handoff" / "Include" localization-status row like the existing ones.

$wb = $excel.ActiveWorkbook

$hyperlinkColor = 15570276   # RGB(100,149,237) == style used by the sheet's existing "HyperLink" cells (FF6495ED)

function Style-AsHyperlink($rng) {
    $rng.Font.Underline = $true
    $rng.Font.Color = $hyperlinkColor
}

function Style-AsDatetime($rng) {
    $rng.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

$fileGuid = "d6da7312-fbe8-491d-afc1-66ef03145004"
$mdName = "$fileGuid.md"
$zhToken = "d7710280bfcd4f242df6d82f942886315d28421d"
$deToken = "d7710280bfcd4f242df6d82f942886315d28421d"
$zhXlfName = "$fileGuid.$zhToken.zh-cn.xlf"
$deXlfName = "$fileGuid.$deToken.de-de.xlf"

$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/9b46e51bf05251965820a1f56f7b1bb6ec630b7f/e2e/$mdName"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f563fa3f04c48794b236d88dbc471217e4f94cc6/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/$zhXlfName"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/36def93dc9e62bd606aa47635dc37430d3e8dca0/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/$deXlfName"

$zhHandoffDatetime = "2016-03-22 11:55:45"
$deHandoffDatetime = "2016-03-22 11:55:53"
$latestHandoffDate = "2016-03-22 11:55:53"
$epoch = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Sheet 1: Overview  (A7:D7)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A7").Value2 = $mdName
Style-AsHyperlink $ws1.Range("A7")
$ws1.Hyperlinks.Add($ws1.Range("A7"), $mdUrl, [System.Type]::Missing, [System.Type]::Missing, $mdName) | Out-Null

$ws1.Range("B7").Value2 = "Ready for handoff"
$ws1.Range("C7").Value2 = "Ready for handoff"

$ws1.Range("D7").Value2 = $latestHandoffDate
Style-AsDatetime $ws1.Range("D7")

# ---------------------------------------------------------------------
# Sheet 2: zh-cn  (A7:L7, sparse like the other data rows)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A7").Value2 = $mdName
Style-AsHyperlink $ws2.Range("A7")
$ws2.Hyperlinks.Add($ws2.Range("A7"), $mdUrl, [System.Type]::Missing, [System.Type]::Missing, $mdName) | Out-Null

$ws2.Range("B7").Value2 = ".md"
$ws2.Range("C7").Value2 = "Ready for handoff"

$ws2.Range("D7").Value2 = $zhXlfName
Style-AsHyperlink $ws2.Range("D7")
$ws2.Hyperlinks.Add($ws2.Range("D7"), $zhXlfUrl, [System.Type]::Missing, [System.Type]::Missing, $zhXlfName) | Out-Null

$ws2.Range("E7").Value2 = $zhHandoffDatetime
Style-AsDatetime $ws2.Range("E7")

$ws2.Range("H7").Value2 = $epoch
Style-AsDatetime $ws2.Range("H7")

$ws2.Range("J7").Value2 = "Include"

# ---------------------------------------------------------------------
# Sheet 3: de-de  (A7:L7, sparse like the other data rows)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A7").Value2 = $mdName
Style-AsHyperlink $ws3.Range("A7")
$ws3.Hyperlinks.Add($ws3.Range("A7"), $mdUrl, [System.Type]::Missing, [System.Type]::Missing, $mdName) | Out-Null

$ws3.Range("B7").Value2 = ".md"
$ws3.Range("C7").Value2 = "Ready for handoff"

$ws3.Range("D7").Value2 = $deXlfName
Style-AsHyperlink $ws3.Range("D7")
$ws3.Hyperlinks.Add($ws3.Range("D7"), $deXlfUrl, [System.Type]::Missing, [System.Type]::Missing, $deXlfName) | Out-Null

$ws3.Range("E7").Value2 = $deHandoffDatetime
Style-AsDatetime $ws3.Range("E7")

$ws3.Range("H7").Value2 = $epoch
Style-AsDatetime $ws3.Range("H7")

$ws3.Range("J7").Value2 = "Include"

Write-Host "Added row 7 ($fileGuid) to Overview, zh-cn and de-de sheets."
